$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting from column Q into the new column R (matches styles
# for the header/spacer row 2 through the data rows 3-6).
$ws.Range("Q2:Q6").Copy()
$ws.Range("R2").PasteSpecial(-4122)

# New 2021 data column
$ws.Range("R3").Value = 2021
$ws.Range("R4").Value = 233306
$ws.Range("R5").Value = 3.5
$ws.Range("R6").Value = 30.8

# Match the selection left behind in the saved worksheet
[void]$ws.Range("Q15").Select()
